$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary values -------------------------------------------------
$ws.Range("E11").Value2 = 805920
$ws.Range("C13").Value2 = 3
$ws.Range("F13").Value2 = 15

# --- Replace the worker/period detail rows (16-31) with the new dataset ------
$data = @(
  @(16, "CC", "73232396",   "EDWIN ALFONSO BUELVAS ARRIETA",   "2306", 37120, 1160000),
  @(17, "CC", "73232396",   "EDWIN ALFONSO BUELVAS ARRIETA",   "2307", 46400, 1160000),
  @(18, "CC", "73232396",   "EDWIN ALFONSO BUELVAS ARRIETA",   "2308", 46400, 1160000),
  @(19, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2409", 52000, 1300000),
  @(20, "CC", "9113176",    "SAUL ENRIQUE COHEN ALVIS",        "2409", 52000, 1300000),
  @(21, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2410", 52000, 1300000),
  @(22, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2411", 52000, 1300000),
  @(23, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2412", 52000, 1300000),
  @(24, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2501", 52000, 1300000),
  @(25, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2502", 52000, 1300000),
  @(26, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2503", 52000, 1300000),
  @(27, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2504", 52000, 1300000),
  @(28, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2505", 52000, 1300000),
  @(29, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2506", 52000, 1300000),
  @(30, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2507", 52000, 1300000),
  @(31, "CC", "1128053543", "JAIRO ANDRES HERNANDEZ BALLESTA", "2508", 52000, 1300000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value2 = $row[1]
    $ws.Range("C$r").Value2 = $row[2]
    $ws.Range("D$r").Value2 = $row[3]
    $ws.Range("E$r").Value2 = $row[4]
    $ws.Range("F$r").Value2 = $row[5]
    $ws.Range("G$r").Value2 = $row[6]
}

# Row 31 becomes the new closing row of the table - copy the formatting that
# the old closing row (43) used (heavier bottom border, etc.) onto it.
$ws.Range("B43:J43").Copy() | Out-Null
$ws.Range("B31:J31").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Drop the now-obsolete trailing detail rows (32-43); this shifts the ----
# --- signature block rows (formerly 48-49) up to become rows 36-37. --------
$ws.Range("32:43").Delete() | Out-Null

Write-Output "done"
